# Auto-generated: refresh market-price derived columns (H-N) per scheduled scrape.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3327.9167
$ws.Range("I32").Value = 2178.2856
$ws.Range("J32").Value = 4937.4
$ws.Range("K32").Value = 2178.2856
$ws.Range("L32").Value = 4937.4
$ws.Range("M32").Value = -1852.2856
$ws.Range("N32").Value = -5589.4
$ws.Range("H132").Value = 3536.7222
$ws.Range("I132").Value = 3536.7222
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10610.1666
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8080.1666
$ws.Range("H138").Value = 1443.5454
$ws.Range("I138").Value = 1443.5454
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 4330.6362
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 809.3638000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1397.1923
$ws.Range("I32").Value = 1253.08
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 1253.08
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -966.0799999999999
$ws.Range("N32").Value = -5574
$ws.Range("H74").Value = 1996
$ws.Range("I74").Value = 1996
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1996
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1122
$ws.Range("H77").Value = 1996
$ws.Range("I77").Value = 1996
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 9980
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -5612

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2218.8823
$ws.Range("I64").Value = 1155.5
$ws.Range("J64").Value = 2546.077
$ws.Range("K64").Value = 1155.5
$ws.Range("L64").Value = 2546.077
$ws.Range("M64").Value = -930.5
$ws.Range("N64").Value = -2996.077
$ws.Range("H67").Value = 2218.8823
$ws.Range("I67").Value = 1155.5
$ws.Range("J67").Value = 2546.077
$ws.Range("K67").Value = 1155.5
$ws.Range("L67").Value = 2546.077
$ws.Range("M67").Value = -375.5
$ws.Range("N67").Value = -4106.077
$ws.Range("H86").Value = 8512.522999999999
$ws.Range("I86").Value = 3523.25
$ws.Range("J86").Value = 15164.889
$ws.Range("K86").Value = 3523.25
$ws.Range("L86").Value = 15164.889
$ws.Range("M86").Value = -2400.25
$ws.Range("N86").Value = -17410.889
$ws.Range("H89").Value = 8512.522999999999
$ws.Range("I89").Value = 3523.25
$ws.Range("J89").Value = 15164.889
$ws.Range("K89").Value = 17616.25
$ws.Range("L89").Value = 75824.44499999999
$ws.Range("M89").Value = -12000.25
$ws.Range("N89").Value = -87056.44499999999
$ws.Range("H105").Value = 3562.6667
$ws.Range("I105").Value = 3688.625
$ws.Range("J105").Value = 2555
$ws.Range("K105").Value = 3688.625
$ws.Range("L105").Value = 2555
$ws.Range("M105").Value = -1941.625
$ws.Range("N105").Value = -6049

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2323.353
$ws.Range("I31").Value = 2080.1
$ws.Range("K31").Value = 2080.1
$ws.Range("M31").Value = -1785.1
$ws.Range("H34").Value = 2323.353
$ws.Range("I34").Value = 2080.1
$ws.Range("K34").Value = 2080.1
$ws.Range("M34").Value = -1878.1
$ws.Range("H47").Value = 12166
$ws.Range("I47").Value = 12166
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 12166
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -11600
$ws.Range("H58").Value = 5649.8335
$ws.Range("I58").Value = 3000
$ws.Range("K58").Value = 3000
$ws.Range("M58").Value = -2797
$ws.Range("H94").Value = 1599.25
$ws.Range("J94").Value = 1599.25
$ws.Range("L94").Value = 1599.25
$ws.Range("N94").Value = -2501.25
$ws.Range("H136").Value = 5649.8335
$ws.Range("I136").Value = 3000
$ws.Range("K136").Value = 9000
$ws.Range("M136").Value = -6450

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 117.25
$ws.Range("I6").Value = 114.333336
$ws.Range("J6").Value = 126
$ws.Range("K6").Value = 343.000008
$ws.Range("L6").Value = 378
$ws.Range("M6").Value = -230.000008
$ws.Range("N6").Value = -604
$ws.Range("H131").Value = 1986.5
$ws.Range("I131").Value = 1986.5
$ws.Range("K131").Value = 5959.5
$ws.Range("M131").Value = -919.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 38098.5
$ws.Range("I3").Value = 1949.5
$ws.Range("J3").Value = 50148.168
$ws.Range("K3").Value = 1949.5
$ws.Range("L3").Value = 50148.168
$ws.Range("M3").Value = -1833.5
$ws.Range("N3").Value = -50380.168
$ws.Range("H10").Value = 533500
$ws.Range("I10").Value = 800000
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 800000
$ws.Range("L10").Value = 500
$ws.Range("M10").Value = -799831
$ws.Range("N10").Value = -838
$ws.Range("H12").Value = 1333333.4
$ws.Range("I12").Value = 1499999.5
$ws.Range("K12").Value = 1499999.5
$ws.Range("M12").Value = -1499859.5
$ws.Range("H70").Value = 8441.083000000001
$ws.Range("I70").Value = 8057.5713
$ws.Range("K70").Value = 8057.5713
$ws.Range("M70").Value = -7787.5713
$ws.Range("H73").Value = 8441.083000000001
$ws.Range("I73").Value = 8057.5713
$ws.Range("K73").Value = 8057.5713
$ws.Range("M73").Value = -7121.5713
$ws.Range("H80").Value = 3299.2
$ws.Range("J80").Value = 4000
$ws.Range("L80").Value = 4000
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 3299.2
$ws.Range("J83").Value = 4000
$ws.Range("L83").Value = 20000
$ws.Range("N83").Value = -29984
$ws.Range("H132").Value = 1622.9286
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 3119.5
$ws.Range("I19").Value = 2740
$ws.Range("K19").Value = 2740
$ws.Range("M19").Value = -2570
$ws.Range("H22").Value = 2597
$ws.Range("I22").Value = 794.5
$ws.Range("K22").Value = 794.5
$ws.Range("M22").Value = -499.5
$ws.Range("H27").Value = 2597
$ws.Range("I27").Value = 794.5
$ws.Range("K27").Value = 794.5
$ws.Range("M27").Value = -687.5
$ws.Range("H46").Value = 2306.8333
$ws.Range("I46").Value = 2231.75
$ws.Range("J46").Value = 2344.375
$ws.Range("K46").Value = 2231.75
$ws.Range("L46").Value = 2344.375
$ws.Range("M46").Value = -2043.75
$ws.Range("N46").Value = -2720.375
$ws.Range("H55").Value = 397.375
$ws.Range("I55").Value = 301
$ws.Range("J55").Value = 493.75
$ws.Range("K55").Value = 301
$ws.Range("L55").Value = 493.75
$ws.Range("M55").Value = -128
$ws.Range("N55").Value = -839.75
$ws.Range("H122").Value = 3263.3076
$ws.Range("I122").Value = 3263.3076
$ws.Range("K122").Value = 9789.9228
$ws.Range("M122").Value = -7339.9228
$ws.Range("H136").Value = 3471.6365
$ws.Range("I136").Value = 1888.2632
$ws.Range("K136").Value = 5664.7896
$ws.Range("M136").Value = -3114.7896

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H81").Value = 4444.5
$ws.Range("I81").Value = 4444.5
$ws.Range("K81").Value = 8889
$ws.Range("M81").Value = -7828
$ws.Range("H84").Value = 4444.5
$ws.Range("I84").Value = 4444.5
$ws.Range("K84").Value = 44445
$ws.Range("M84").Value = -39141
$ws.Range("H122").Value = 3078.3572
$ws.Range("I122").Value = 2887.2
$ws.Range("J122").Value = 4671.3335
$ws.Range("K122").Value = 8661.599999999999
$ws.Range("L122").Value = 14014.0005
$ws.Range("M122").Value = -6211.599999999999
$ws.Range("N122").Value = -18914.0005
